$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null -and $val -is [string]) {
            if ($val -match "D51|D64|D80|S30") {
                $newVal = $val -replace "D64", "D69"
                $newVal = $newVal -replace "D51", "D55"
                $newVal = $newVal -replace "D80", "D86"
                $newVal = $newVal -replace "S30", "S31"
                $cell.Value2 = $newVal
            }
        }
    }
}
